$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace the old header row with a single title cell.
$ws.Range("B1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# Drop the trailing "Level" / "Course" columns entirely (K:L),
# which also shrinks the used range down to column J.
$ws.Columns("K:L").Delete()

# Re-key each data row (2-10) into the new column layout:
#   A Day | B Time | C Hours | D Module Code | E Module Title |
#   F Class Type | G Lecturer | H Group | I Block | J Room
$ws.Range("A2").Value = "SUN"
$ws.Range("B2").Value = "12:30-15:00"
$ws.Range("C2").Value = 2.5
$ws.Range("D2").Value = "5CS020"
$ws.Range("E2").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F2").Value = "Workshop"
$ws.Range("G2").Value = "Mr. Prabin Sapkota"
$ws.Range("H2").Value = "L5CG2"
$ws.Range("I2").Value = "HCK"
$ws.Range("J2").Value = "Lab-03 Gahanapokhari"

$ws.Range("A3").Value = "MON"
$ws.Range("B3").Value = "9:30-12:00"
$ws.Range("C3").Value = 2.5
$ws.Range("D3").Value = "5CS024"
$ws.Range("E3").Value = "Collaborative Development"
$ws.Range("F3").Value = "Workshop"
$ws.Range("G3").Value = "Mr. Anmol Adhikari"
$ws.Range("H3").Value = "L5CG2"
$ws.Range("I3").Value = "WLV"
$ws.Range("J3").Value = "TR-01 Dudley"

$ws.Range("A4").Value = "TUE"
$ws.Range("B4").Value = "9:00-11:00"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "5CS022"
$ws.Range("E4").Value = "Human Computer Interaction"
$ws.Range("F4").Value = "Lecture"
$ws.Range("G4").Value = "Mr. Pravash Karki"
$ws.Range("H4").Value = "L5CG(1+2+3+4)"
$ws.Range("I4").Value = "WLV"
$ws.Range("J4").Value = "LT-02 Telford"

$ws.Range("A5").Value = "TUE"
$ws.Range("B5").Value = "12:00-14:00"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "5CS020"
$ws.Range("E5").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F5").Value = "Lecture"
$ws.Range("G5").Value = "Mr. Sumanta Silwal"
$ws.Range("H5").Value = "L5CG(1+2+3+4)"
$ws.Range("I5").Value = "WLV"
$ws.Range("J5").Value = "LT-01 Wulfruna"

$ws.Range("A6").Value = "WED"
$ws.Range("B6").Value = "9:30-11:30"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "5CS024"
$ws.Range("E6").Value = "Collaborative Development"
$ws.Range("F6").Value = "Lecture"
$ws.Range("G6").Value = "Mr. Udaya Kandel"
$ws.Range("H6").Value = "L5CG(1+2+3+4)"
$ws.Range("I6").Value = "WLV"
$ws.Range("J6").Value = "LT-01 Wulfruna"

$ws.Range("A7").Value = "WED"
$ws.Range("B7").Value = "12:30-14:30"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "5CS020"
$ws.Range("E7").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F7").Value = "Tutorial"
$ws.Range("G7").Value = "Mr. Prabin Sapkota"
$ws.Range("H7").Value = "L5CG2"
$ws.Range("I7").Value = "WLV"
$ws.Range("J7").Value = "TR-03 Westbromwich"

$ws.Range("A8").Value = "THU"
$ws.Range("B8").Value = "9:00-11:00"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "5CS022"
$ws.Range("E8").Value = "Human Computer Interaction"
$ws.Range("F8").Value = "Tutorial"
$ws.Range("G8").Value = "Mr. Dipesh Shrestha"
$ws.Range("H8").Value = "L5CG2"
$ws.Range("I8").Value = "WLV"
$ws.Range("J8").Value = "TR-02 Stafford"

$ws.Range("A9").Value = "FRI"
$ws.Range("B9").Value = "10:00-12:00"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "5CS024"
$ws.Range("E9").Value = "Collaborative Development"
$ws.Range("F9").Value = "Tutorial"
$ws.Range("G9").Value = "Mr. Anmol Adhikari"
$ws.Range("H9").Value = "L5CG2"
$ws.Range("I9").Value = "HCK"
$ws.Range("J9").Value = "TR-09  Chandragiri"

$ws.Range("A10").Value = "FRI"
$ws.Range("B10").Value = "13:30-16:00"
$ws.Range("C10").Value = 2.5
$ws.Range("D10").Value = "5CS022"
$ws.Range("E10").Value = "Human Computer Interaction"
$ws.Range("F10").Value = "Workshop"
$ws.Range("G10").Value = "Mr. Dipesh Shrestha"
$ws.Range("H10").Value = "L5CG2"
$ws.Range("I10").Value = "WLV"
$ws.Range("J10").Value = "SR-03 Wolves"
